# Update "想去人数" (want-to-go count) figures that were refreshed when the
# site's data was regenerated (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 59
$ws1.Range("F3").Value  = 11591
$ws1.Range("F4").Value  = 212
$ws1.Range("F5").Value  = 330
$ws1.Range("F7").Value  = 11556
$ws1.Range("F8").Value  = 478
$ws1.Range("F11").Value = 1758
$ws1.Range("F12").Value = 5728
$ws1.Range("F13").Value = 113

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 4

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 59
$ws4.Range("F5").Value  = 11591
$ws4.Range("F6").Value  = 212
$ws4.Range("F7").Value  = 330
$ws4.Range("F9").Value  = 11556
$ws4.Range("F10").Value = 478
$ws4.Range("F13").Value = 1758
$ws4.Range("F14").Value = 4
$ws4.Range("F15").Value = 5728
$ws4.Range("F16").Value = 113
